$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): column F "想去人数" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2665
$ws1.Range("F3").Value = 578
$ws1.Range("F4").Value = 471
$ws1.Range("F5").Value = 301
$ws1.Range("F6").Value = 204
$ws1.Range("F7").Value = 491
$ws1.Range("F8").Value = 1231
$ws1.Range("F9").Value = 574
$ws1.Range("F12").Value = 130
$ws1.Range("F13").Value = 365
$ws1.Range("F14").Value = 5766
$ws1.Range("F15").Value = 87
$ws1.Range("F16").Value = 1794
$ws1.Range("F17").Value = 4204
$ws1.Range("F18").Value = 438
$ws1.Range("F19").Value = 240
$ws1.Range("F20").Value = 304
$ws1.Range("F21").Value = 4911
$ws1.Range("F22").Value = 6285
$ws1.Range("F26").Value = 3794
$ws1.Range("F27").Value = 503
$ws1.Range("F28").Value = 69
$ws1.Range("F31").Value = 996
$ws1.Range("F32").Value = 1420
$ws1.Range("F33").Value = 482
$ws1.Range("F34").Value = 572
$ws1.Range("F35").Value = 1613
$ws1.Range("F36").Value = 206
$ws1.Range("F37").Value = 1737
$ws1.Range("F38").Value = 205
$ws1.Range("F39").Value = 1149
$ws1.Range("F40").Value = 1338
$ws1.Range("F41").Value = 637
$ws1.Range("F43").Value = 3440
$ws1.Range("F45").Value = 295
$ws1.Range("F48").Value = 20
$ws1.Range("F49").Value = 3900

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 1212

# Sheet "本地生活" (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 3964

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3964
$ws4.Range("F4").Value = 578
$ws4.Range("F5").Value = 471
$ws4.Range("F6").Value = 301
$ws4.Range("F7").Value = 1212
$ws4.Range("F10").Value = 204
$ws4.Range("F11").Value = 491
$ws4.Range("F13").Value = 1231
$ws4.Range("F14").Value = 574
$ws4.Range("F16").Value = 130
$ws4.Range("F17").Value = 365
$ws4.Range("F18").Value = 1794
$ws4.Range("F19").Value = 4204
$ws4.Range("F20").Value = 4912
$ws4.Range("F22").Value = 1061
$ws4.Range("F24").Value = 3794
$ws4.Range("F25").Value = 503
$ws4.Range("F26").Value = 69
$ws4.Range("F29").Value = 1420
$ws4.Range("F30").Value = 482
$ws4.Range("F31").Value = 572
$ws4.Range("F32").Value = 1613
$ws4.Range("F33").Value = 206
$ws4.Range("F34").Value = 1737
$ws4.Range("F37").Value = 637
$ws4.Range("F41").Value = 3440
$ws4.Range("F44").Value = 295
$ws4.Range("F48").Value = 3900
